$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.75
$ws.Range("I2").Value = 2.38
$ws.Range("J2").Value = 4.75
$ws.Range("K2").Value = 1.73
$ws.Range("L2").Value = 3.4
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("S2").Value = 1.8
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 2.63
$ws.Range("V2").Value = 1.44
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 17
$ws.Range("AB2").Value = 67
$ws.Range("AD2").Value = 6
$ws.Range("AG2").Value = 5
$ws.Range("AH2").Value = 9.5
$ws.Range("AJ2").Value = 23
$ws.Range("AN2").Value = 5
$ws.Range("AO2").Value = 26
$ws.Range("AQ2").Value = 101
$ws.Range("AT2").Value = 1.91
$ws.Range("AW2").Value = 4
$ws.Range("AX2").Value = 17
$ws.Range("AZ2").Value = 51

# Row 3 updates
$ws.Range("G3").Value = 1.62
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("AG3").Value = 13
$ws.Range("AK3").Value = 51
$ws.Range("AN3").Value = 3.4
$ws.Range("AU3").Value = 9.5
